$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2" = "59.249.96"
    "E2" = "  -6.05%  "
    "D3" = "2.451.03"
    "E3" = "  -8.87%  "
    "D4" = "0.999"
    "E4" = "  -0.15%  "
    "D5" = "540.01"
    "E5" = "  -3.02%  "
    "D6" = "147.44"
    "E6" = "  -7.44%  "
    "D8" = "0.570"
    "E8" = "  -3.82%  "
    "D9" = "2.465.99"
    "E9" = "  -8.58%  "
    "D10" = "0.0994"
    "E10" = "  -6.61%  "
    "D12" = "5.33"
    "E12" = "  -1.46%  "
    "D13" = "0.353"
    "E13" = "  -5.14%  "
    "D14" = "2.891.98"
    "E14" = "  -8.67%  "
    "D15" = "24.09"
    "E15" = "  -9.53%  "
    "D16" = "59.046.43"
    "E16" = "  -6.25%  "
    "E17" = "  -6.69%  "
    "D18" = "2.523.37"
    "E18" = "  -6.24%  "
    "D19" = "11.16"
    "E19" = "  -7.03%  "
    "D20" = "4.36"
    "E20" = "  -6.11%  "
    "D21" = "325.03"
    "E21" = "  -6.28%  "
    "D22" = "0.966"
    "E22" = "  -3.46%  "
    "D23" = "5.74"
    "E23" = "  -9.32%  "
    "D24" = "0.460"
    "E24" = "  -10.22%  "
    "D25" = "60.76"
    "E25" = "  -4.42%  "
    "E26" = "  -4.81%  "
    "D27" = "0.980"
    "E27" = "  -1.83%  "
    "D28" = "7.72"
    "E28" = "  -6.71%  "
    "D29" = "1.29"
    "E29" = "  -11.47%  "
    "E30" = "  -6.58%  "
    "B31" = "PEPE"
    "C31" = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
    "D31" = "0.0₃0773"
    "E31" = "  -10.64%  "
    "B32" = "Aptos"
    "C32" = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
    "D32" = "6.68"
    "E32" = "  -8.75%  "
    "D34" = "157.18"
    "E34" = "  -4.44%  "
    "D35" = "4.52"
    "E35" = "  -8.95%  "
    "E36" = "  -8.58%  "
    "D37" = "18.42"
    "E37" = "  -5.95%  "
    "E38" = "  -4.23%  "
    "B39" = "Bittensor"
    "C39" = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
    "D39" = "318.85"
    "E39" = "  -11.97%  "
    "B40" = "RenderToken"
    "C40" = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
    "D40" = "5.88"
    "E40" = "  -10.34%  "
    "D41" = "36.57"
    "E41" = "  -5.21%  "
    "D42" = "0.833"
    "E42" = "  -13.91%  "
    "D43" = "3.69"
    "E43" = "  -8.61%  "
    "E44" = "  -0.23%  "
    "D45" = "10.72"
    "E45" = "  -3.01%  "
    "D46" = "0.0944"
    "E46" = "  -3.26%  "
    "D47" = "0.583"
    "E47" = "  -6.21%  "
    "D48" = "0.0526"
    "E48" = "  -7.06%  "
    "D49" = "19.03"
    "E49" = "  -10.32%  "
    "D50" = "0.0229"
    "E50" = "  -6.32%  "
    "B51" = "Aave"
    "C51" = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
    "D51" = "121.69"
    "E51" = "  -6.32%  "
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}

